{"js": "// The document has three headers (Primary, FirstPage, EvenPages -- i.e.\n// header2.xml, header3.xml, header1.xml in the underlying package) that each\n// contain a legacy VML watermark picture (\"WordPictureWatermark...\" shape,\n// rendered through a <w:pict> fallback run). The edit removes that\n// watermark shape/run from every header while leaving the rest of each\n// header (the \"DRAFT VERSION - COMMENTS WELCOME\" text/hyperlink, the tab\n// runs, and the OWASP logo picture) untouched.\n\ncontext.document.sections.load(\"items\");\nawait context.sync();\n\nconst headerTypes = [\"Primary\", \"FirstPage\", \"EvenPages\"];\n\n// Collect every watermark-like floating shape across all three header\n// stories first, then delete them. (Each header only ever has the one\n// watermark shape, but we match by name defensively in case of stray\n// shapes.)\nconst shapesToDelete = [];\nfor (const type of headerTypes) {\n  const header = context.document.sections.items[0].getHeader(type);\n  const shapes = header.shapes;\n  shapes.load(\"items/name\");\n  await context.sync();\n\n  for (const shape of shapes.items) {\n    if (shape.name && shape.name.indexOf(\"WordPictureWatermark\") === 0) {\n      shapesToDelete.push(shape);\n    }\n  }\n}\n\nfor (const shape of shapesToDelete) {\n  shape.delete();\n}\n\nawait context.sync();\n", "ps1": "# The document's three headers (Primary, First Page, Even Pages -- i.e.\n# header2.xml, header3.xml, header1.xml in the underlying package) each\n# contain a legacy VML watermark picture (\"WordPictureWatermark...\"\n# shape, surfaced through a <w:pict> fallback run). Remove that\n# watermark shape from every header/section while leaving the rest of\n# each header (the \"DRAFT VERSION - COMMENTS WELCOME\" text/hyperlink,\n# the tab runs, and the OWASP logo picture) untouched.\n\n$d = $word.ActiveDocument\n\nforeach ($sec in $d.Sections) {\n    for ($i = 1; $i -le 3; $i++) {\n        $hdr = $sec.Headers.Item($i)\n        if (-not $hdr.Exists) { continue }\n\n        for ($j = $hdr.Shapes.Count; $j -ge 1; $j--) {\n            $shp = $hdr.Shapes.Item($j)\n            if ($shp.Name -like \"WordPictureWatermark*\") {\n                $shp.Delete()\n            }\n        }\n    }\n}\n"}
